$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (RF)
$ws.Range("B3").Value = 0.898
$ws.Range("C3").Value = 0.89
$ws.Range("D3").Value = 0.112
$ws.Range("E3").Value = 0.335
$ws.Range("F3").Value = 0.241
$ws.Range("G3").Value = 0.972

# Row 4 (NN)
$ws.Range("B4").Value = 0.803
$ws.Range("C4").Value = 0.787
$ws.Range("D4").Value = 0.216
$ws.Range("E4").Value = 0.465
$ws.Range("F4").Value = 0.349
$ws.Range("G4").Value = 0.936

# Row 5 (RNN)
$ws.Range("B5").Value = 0.613
$ws.Range("C5").Value = 0.596
$ws.Range("D5").Value = 0.422
$ws.Range("E5").Value = 0.65
$ws.Range("F5").Value = 0.468
$ws.Range("G5").Value = 0.784
